$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add new "Save" column header in H1, matching the style of the existing
# header cells (e.g. G1 "sum") by copying the cell formatting over.
$ws.Range("G1").Copy($ws.Range("H1"))
$ws.Range("H1").Value = "Save"

# Populate the new Save column (H2:H16) with its values.
$saveValues = @{
    2  = 0
    3  = 0
    4  = 0
    5  = 1
    6  = 0
    7  = 1
    8  = 1
    9  = 1
    10 = 0
    11 = 1
    12 = 1
    13 = 0
    14 = 0
    15 = 1
    16 = 1
}

foreach ($row in $saveValues.Keys) {
    $ws.Cells.Item($row, 8).Value = $saveValues[$row]
}
